$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all data rows (2-51)
# from 45186 to 45188
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 3).Value = 45188
}

# Row 51 gains an explicit row height (15, customHeight)
$ws.Rows.Item(51).RowHeight = 15

# Add the new row 52 with the new case data
$ws.Range("A52").Value = "A 43792-2023"
$ws.Range("B52").Value = 45187
$ws.Range("B52").NumberFormat = "YYYY-MM-DD"
$ws.Range("C52").Value = 45188
$ws.Range("C52").NumberFormat = "YYYY-MM-DD"
$ws.Range("D52").Value = "SKÅNE LÄN"
$ws.Range("E52").Value = "SIMRISHAMN"
$ws.Range("G52").Value = 0.7
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = 0
$ws.Range("N52").Value = 0
$ws.Range("O52").Value = 0
$ws.Range("P52").Value = 0
$ws.Range("Q52").Value = 0
$ws.Range("R52").Value = ""
$ws.Range("R52").WrapText = $true
